$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" (totals) sheet: insert a new "2022-Q4" row at the top
#    of the data (row 2), shifting the existing rows down by one and
#    appending the previously-last row ("2021-Q1") as the new row 8.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Shift existing data (rows 2..7) down to rows 3..8, bottom-up so we never
# clobber a row before it has been read.
for ($r = 7; $r -ge 2; $r--) {
    $dest = $r + 1
    $total.Cells.Item($dest, 2).Value = $total.Cells.Item($r, 2).Value()
    $total.Cells.Item($dest, 3).Value = $total.Cells.Item($r, 3).Value()
    $total.Cells.Item($dest, 4).Value = $total.Cells.Item($r, 4).Value()
}

# The new last row (row 8) needs the same formatting as the other index
# cells in column A (style copied from row 7, which already carries it).
$total.Range("A7").Copy()
$total.Range("A8").PasteSpecial(-4122)  # xlPasteFormats

# Write the brand-new "2022-Q4" row.
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 1
$total.Cells.Item(2, 4).Value = 0.01

# Re-number the index column (A) 0..6 for rows 2..8.
for ($r = 2; $r -le 8; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------------
# 2. Insert a brand-new "2022-Q4" worksheet (holding fund-holding detail)
#    right after "总计". Clone the "2022-Q1" sheet (which already has the
#    exact same 1-header-row + 1-data-row shape and formatting we need) so
#    all styling/borders/page setup carry over, then overwrite its values.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q1")
$template.Copy($null, $total)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Header row (row 1) - text labels, identical to the other quarter sheets.
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Data row (row 2).
$q4.Range("A2").Value = 0

$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "165531"

$q4.Range("C2").NumberFormat = "@"
$q4.Range("C2").Value = "信诚多策略灵活配置混合（LOF）"

$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "0.89"

$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "72.25"

$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "1.43"

$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.0127"

$q4.Range("H2").Value = 1

# ---------------------------------------------------------------------------
# 3. The sheet still named "2021-Q3" gets its "基金金额" column header
#    relabelled to "基金规模" (matching the other quarters).
# ---------------------------------------------------------------------------
$q3_2021 = $wb.Worksheets.Item("2021-Q3")
$q3_2021.Range("D1").Value = "基金规模"
